$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "70-8=",
    "45-35=",
    "32+17=",
    "55+0=",
    "97-32=",
    "45+19=",
    "30+35=",
    "76-13=",
    "15+66=",
    "8+14=",
    "80-74=",
    "99-51=",
    "38+46=",
    "81-71=",
    "17+41=",
    "26+19=",
    "61+32=",
    "63+25=",
    "43-26=",
    "44+41=",
    "20+34=",
    "19+9=",
    "49+47=",
    "31+20=",
    "15+61=",
    "92-47=",
    "57+40=",
    "43-20=",
    "3+43=",
    "19+49=",
    "4+5=",
    "1+84=",
    "44-42=",
    "35-28=",
    "39+49=",
    "6-2=",
    "31+34=",
    "13+69=",
    "49-26=",
    "35-32=",
    "64-55=",
    "96-96=",
    "52-28=",
    "17+32=",
    "90-38=",
    "27+51=",
    "95-70=",
    "27+2=",
    "92+6=",
    "29+11=",
    "93+6=",
    "75+14=",
    "80-7=",
    "83-67=",
    "29-10=",
    "69+8=",
    "16+28=",
    "67-0=",
    "3+46=",
    "36+13=",
    "45+19=",
    "67-66=",
    "28+52=",
    "11+40=",
    "87-72=",
    "46-24=",
    "31-5=",
    "71-23=",
    "11-8=",
    "97-36=",
    "43-25=",
    "10+88=",
    "40+9=",
    "14+54=",
    "64-46=",
    "61-50=",
    "96-36=",
    "35-29=",
    "41-9=",
    "41-7=",
    "30-5=",
    "16+4=",
    "41+43=",
    "3+15=",
    "47+39=",
    "51-36=",
    "42+20=",
    "49+32=",
    "21-0=",
    "93-63=",
    "67-5=",
    "25-17=",
    "4+80=",
    "3+62=",
    "46-23=",
    "21+64=",
    "1+19=",
    "36+5=",
    "36-1=",
    "13+39="
)

$numCols = 5
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated $idx cells"
